$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date/time demo values in B3 (date) and C3 (time)
$ws.Range("B3").Value = 44561
$ws.Range("C3").Value = 0.2488078703703704
